$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded at the top of the Mango
# data block (row 46). Insert a fresh row there, which pushes the
# pre-existing rows 46:70 down to 47:71 (the former last row, 44595,
# ends up at row 71).
$ws.Rows(46).Insert()

# Populate the new row 46 with the new data point.
$ws.Range("A46").Value = 7
$ws.Range("B46").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C46").Value = 'Ñuble'
$ws.Range("D46").Value = 44596
$ws.Range("E46").Value = 16
$ws.Range("F46").Value = 'Fruta'
$ws.Range("G46").Value = 100108
$ws.Range("H46").Value = 'Tropicales y subtropicales'
$ws.Range("I46").Value = 100108002
$ws.Range("J46").Value = 'Mango'
$ws.Range("K46").Value = 'Sin especificar'
$ws.Range("L46").Value = 'Primera'
$ws.Range("M46").Value = 100
$ws.Range("N46").Value = 7000
$ws.Range("O46").Value = 7500
$ws.Range("P46").Value = 7250
$ws.Range("Q46").Value = '$/bandeja 4 kilos'
$ws.Range("R46").Value = 'Perú'
$ws.Range("S46").Value = 1812
$ws.Range("T46").Value = 4
